# Fix #8177 / #8211: Navigator download should include internationalized
# label-xx / description-xx columns on the "entities" and "attributes"
# sheets of the generated EMX workbook, and should not write out empty
# expression values.
#
# This script reproduces the resulting workbook layout:
#  - "entities" sheet (sheet 2) gains 16 new header columns (I..X) for
#    label-en/description-en .. label-xx/description-xx, and the existing
#    (always-empty) "description" / "extends" data cells are cleared.
#  - "attributes" sheet (sheet 3) gains the same 16 new header columns
#    (X..AM) for label-en/description-en .. label-xx/description-xx.
#  - the "attributes" sheet becomes the active/selected tab (it was the
#    "packages" sheet before).

$wb = $excel.ActiveWorkbook

$newHeaders = @("label-en","description-en","label-nl","description-nl","label-de","description-de","label-es","description-es","label-it","description-it","label-pt","description-pt","label-fr","description-fr","label-xx","description-xx")

# ---------------------------------------------------------------------
# entities sheet: add the 16 new header columns after "tags" (column H)
# ---------------------------------------------------------------------
$entities = $wb.Worksheets.Item("entities")

$col = 9
foreach ($h in $newHeaders) {
    $entities.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}

# the "description" (D) and "extends" (F) values on the data rows were
# already empty strings; the Navigator download no longer emits them
$entities.Cells.Item(2, 4).Value = ""
$entities.Cells.Item(3, 4).Value = ""
$entities.Cells.Item(4, 4).Value = ""
$entities.Cells.Item(2, 6).Value = ""
$entities.Cells.Item(3, 6).Value = ""
$entities.Cells.Item(4, 6).Value = ""

# ---------------------------------------------------------------------
# attributes sheet: add the 16 new header columns after "tags" (column W)
# ---------------------------------------------------------------------
$attributes = $wb.Worksheets.Item("attributes")

$col = 24
foreach ($h in $newHeaders) {
    $attributes.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}

# ---------------------------------------------------------------------
# the attributes sheet is now the active tab (was "packages" before)
# ---------------------------------------------------------------------
$attributes.Activate()
